$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new line to the "Switching State" notes in B5, and grow the row
# to fit the extra line of wrapped text. (Do this first so the updated
# shared string is appended to the table before the new row's strings.)
$switchingStateNotes = $ws.Range("B5").Value()
$ws.Range("B5").Value = $switchingStateNotes + "`nExtract events + their data from the coordinator and store somewhere else between states."
$ws.Rows.Item(5).RowHeight = 105

# Insert a new row at position 6 (pushes existing rows 6-17 down to 7-18),
# copying formatting (e.g. wrap-text style) from the row above (row 5).
$ws.Rows.Item(6).Insert()

# New row 6: "Multithreading" task.
$ws.Range("A6").Value = "Multithreading"
$ws.Range("B6").Value = "Render system and camera system need multithreading "

# Update the view: selection moved to A8 (this also clears the old
# scrolled-down topLeftCell="A9" setting).
[void]$ws.Range("A8").Select()
